$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names (shared strings used by column B)
$ws.Range("B2").Value = "John Doe"
$ws.Range("B3").Value = "Jane Doe"

# Update numeric values
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 19
$ws.Range("A3").Value = 5
